$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.27 = 8667.88 pesos`n✅ 8667.88 pesos = 2.26 = 956.64 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 439.9
$wsTasas.Range("O10").Value = 3813
$wsTasas.Range("N12").Value = 3830
$wsTasas.Range("O12").Value = 422.701
